# powerblade_3v2.xlsx - add per-part "qty per board" column (L) while
# assembling pogos.
#
# Column L gets a header-less numeric/text quantity value for every part
# row (2-62) except row 35 (no value there), right-aligned to match the
# "10+4" style entries (special two-part quantity notations).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlRight = -4152
$xlRight = -4152

$colLValues = @{
    2=10; 3=5; 4=4; 5=4; 6=8; 7=5; 8=5; 9=5; 10="10+4"; 11=10;
    12=10; 13=5; 14=50; 15=10; 16=10; 17=4; 18=10; 19=15; 20=10; 21=10;
    22=4; 23=4; 24=3; 25=10; 26=5; 27=10; 28=7; 29=4; 30=10; 31="10+4";
    32=4; 33=10; 34=10; 36="NP"; 37="NP"; 38="NP"; 39=5; 40=50; 41=20;
    42=4; 43=20; 44=20; 45=20; 46=20; 47=20; 48=50; 49=20; 50=20; 51=50;
    52=20; 53=20; 54=20; 55=4; 56=20; 57=50; 58=20; 59=10; 60=10; 61=10;
    62=10
}

foreach ($row in $colLValues.Keys) {
    $cell = $ws.Range("L$row")
    $cell.Value = $colLValues[$row]
    $cell.HorizontalAlignment = $xlRight
}

# Move the selection/viewport to where editing left off.
$ws.Range("A41").Select()
